# Add a new "Address" column (F) before the existing "District" column,
# which shifts District from F to G. Populate the new Address column by
# extracting the school/institution portion of the multi-line "NAMES"
# text already present in column B (everything on the second line of the
# cell except for the last comma-separated segment, which is the
# district and is already captured in the District column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; this pushes the existing "District"
# column (and all its data) from F to G, preserving values/styles.
$ws.Columns("F").Insert()

# New header for the inserted column.
$ws.Range("F2").Value = "Address"

$lastRow = 74

for ($r = 3; $r -le $lastRow; $r++) {
    $text = $ws.Range("B$r").Value2
    if ($null -eq $text -or $text -eq "") {
        continue
    }

    $lines = $text -split "`n"
    if ($lines.Length -ge 2) {
        $addrLine = $lines[1]
    } else {
        $addrLine = $lines[0]
    }

    $parts = $addrLine -split ","
    if ($parts.Length -le 1) {
        # No comma-separated district segment to split off, so there is
        # nothing meaningful to place in the Address column (the
        # existing District value already moved to column G via the
        # column insert above).
        continue
    }

    $n = $parts.Length
    $addrParts = $parts[0..($n - 2)]
    $addr = ""
    foreach ($p in $addrParts) {
        $addr = $addr + $p.Trim()
    }

    $ws.Range("F$r").Value = $addr
}
